$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder columns: Mollie-ID, ID-nummer, Naam, Aantal
$ws.Range("A1").Value = "Mollie-ID"
$ws.Range("B1").Value = "ID-nummer"
$ws.Range("C1").Value = "Naam"
$ws.Range("D1").Value = "Aantal"

$ws.Range("A2").Value = "mollie_ABC123"
$ws.Range("B2").Value = 1445758
$ws.Range("C2").Value = "Ben Gortemaker"
$ws.Range("D2").Value = 1
